$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.687074477498015
$ws.Range("C2").Value = 1.962075346092774
$ws.Range("D2").Value = 0.6917674849873663
$ws.Range("E2").Value = 0.2813313834455471
$ws.Range("G2").Value = 0.00265223839910164
$ws.Range("I2").Value = 5.249643776912166
$ws.Range("J2").Value = 0.1457323425992385
$ws.Range("B3").Value = 2.629574434925587
$ws.Range("C3").Value = 1.879422835144055
$ws.Range("D3").Value = 0.6841789757318395
$ws.Range("E3").Value = 0.2774161609282473
$ws.Range("G3").Value = 0.002661996347303691
$ws.Range("I3").Value = 4.933919875177537
$ws.Range("J3").Value = 0.1430009961348091
$ws.Range("B4").Value = 2.597605440569851
$ws.Range("C4").Value = 1.830172562647363
$ws.Range("D4").Value = 0.6799450801667035
$ws.Range("E4").Value = 0.2751937454154643
$ws.Range("G4").Value = 0.002668277244588583
$ws.Range("I4").Value = 4.740498520279317
$ws.Range("J4").Value = 0.1414242038296862
$ws.Range("B5").Value = 2.585407838365768
$ws.Range("C5").Value = 1.810472964790108
$ws.Range("D5").Value = 0.6783258511028407
$ws.Range("E5").Value = 0.2743333149010283
$ws.Range("G5").Value = 0.002670909927310154
$ws.Range("I5").Value = 4.661766826334969
$ws.Range("J5").Value = 0.1408065748754765
$ws.Range("B6").Value = 2.583432311785543
$ws.Range("C6").Value = 1.807224037073468
$ws.Range("D6").Value = 0.6780633636024334
$ws.Range("E6").Value = 0.2741931600366598
$ws.Range("G6").Value = 0.002671351512296284
$ws.Range("I6").Value = 4.648698286085249
$ws.Range("J6").Value = 0.1407055153614181
$ws.Range("B7").Value = 2.597437590268839
$ws.Range("C7").Value = 1.829905396260529
$ws.Range("D7").Value = 0.6799228141412925
$ws.Range("E7").Value = 0.2751819588046516
$ws.Range("G7").Value = 0.002668312453121553
$ws.Range("I7").Value = 4.739436382846463
$ws.Range("J7").Value = 0.1414157737138737
$ws.Range("B8").Value = 2.666550549870067
$ws.Range("C8").Value = 1.933261670135948
$ws.Range("D8").Value = 0.689062121525609
$ws.Range("E8").Value = 0.2799434705092239
$ws.Range("G8").Value = 0.002655543085293884
$ws.Range("I8").Value = 5.14067969973209
$ws.Range("J8").Value = 0.1447695887030775
$ws.Range("B9").Value = 2.828979289754955
$ws.Range("C9").Value = 2.148150250608637
$ws.Range("D9").Value = 0.7104025685030706
$ws.Range("E9").Value = 0.2907420157758622
$ws.Range("G9").Value = 0.002632782043708433
$ws.Range("I9").Value = 5.931919908146625
$ws.Range("J9").Value = 0.1521557101994659
$ws.Range("B10").Value = 2.965341798777558
$ws.Range("C10").Value = 2.31393649761867
$ws.Range("D10").Value = 0.7282282407530545
$ws.Range("E10").Value = 0.2995972713864674
$ws.Range("G10").Value = 0.002617424905153104
$ws.Range("I10").Value = 6.51732333652464
$ws.Range("J10").Value = 0.1580959511421725
$ws.Range("B11").Value = 3.031217255711795
$ws.Range("C11").Value = 2.391179793726337
$ws.Range("D11").Value = 0.7368181932760081
$ws.Range("E11").Value = 0.3038328931438912
$ws.Range("G11").Value = 0.002610729731212034
$ws.Range("I11").Value = 6.784837002094548
$ws.Range("J11").Value = 0.1609144779455391
$ws.Range("B12").Value = 3.056727193443123
$ws.Range("C12").Value = 2.420701063684078
$ws.Range("D12").Value = 0.740141331128342
$ws.Range("E12").Value = 0.3054672022605089
$ws.Range("G12").Value = 0.00260823585004453
$ws.Range("I12").Value = 6.886336772260506
$ws.Range("J12").Value = 0.1619988810747657
$ws.Range("B13").Value = 3.051207882967617
$ws.Range("C13").Value = 2.414330946199925
$ws.Range("D13").Value = 0.7394224881789171
$ws.Range("E13").Value = 0.3051138646761018
$ws.Range("G13").Value = 0.002608771115411912
$ws.Range("I13").Value = 6.864467778093342
$ws.Range("J13").Value = 0.1617645694725098
$ws.Range("B14").Value = 3.033304595331401
$ws.Range("C14").Value = 2.39360304843683
$ws.Range("D14").Value = 0.7370901739961084
$ws.Range("E14").Value = 0.3039667367834866
$ws.Range("G14").Value = 0.002610523729888126
$ws.Range("I14").Value = 6.79318333799273
$ws.Range("J14").Value = 0.1610033478212074
$ws.Range("B15").Value = 3.022412160493502
$ws.Range("C15").Value = 2.380942143871607
$ws.Range("D15").Value = 0.7356707553477122
$ws.Range("E15").Value = 0.3032680599653617
$ws.Range("G15").Value = 0.002611602641755897
$ws.Range("I15").Value = 6.749546089892874
$ws.Range("J15").Value = 0.1605393140851703
$ws.Range("B16").Value = 2.961115008717911
$ws.Range("C16").Value = 2.308925918766363
$ws.Range("D16").Value = 0.7276766424676282
$ws.Range("E16").Value = 0.2993246826604619
$ws.Range("G16").Value = 0.002617868270331859
$ws.Range("I16").Value = 6.49986709365615
$ws.Range("J16").Value = 0.1579141247506044
$ws.Range("B17").Value = 2.924503744792332
$ws.Range("C17").Value = 2.265219868339386
$ws.Range("D17").Value = 0.7228965023522278
$ws.Range("E17").Value = 0.2969590394067652
$ws.Range("G17").Value = 0.002621786251585474
$ws.Range("I17").Value = 6.347024124595976
$ws.Range("J17").Value = 0.1563336957063086
$ws.Range("B18").Value = 2.903806956813241
$ws.Range("C18").Value = 2.240252691736373
$ws.Range("D18").Value = 0.7201922916694343
$ws.Range("E18").Value = 0.2956178634368101
$ws.Range("G18").Value = 0.00262406717182136
$ws.Range("I18").Value = 6.259225373029011
$ws.Range("J18").Value = 0.155435596948891
$ws.Range("B19").Value = 2.896861041607508
$ws.Range("C19").Value = 2.231828422655326
$ws.Range("D19").Value = 0.7192844212064244
$ws.Range("E19").Value = 0.2951670932311004
$ws.Range("G19").Value = 0.00262484417029733
$ws.Range("I19").Value = 6.229516820578681
$ws.Range("J19").Value = 0.1551333810715505
$ws.Range("B20").Value = 2.92836362683056
$ws.Range("C20").Value = 2.269854654879964
$ws.Range("D20").Value = 0.7234006702827571
$ws.Range("E20").Value = 0.2972088461387159
$ws.Range("G20").Value = 0.002621366342699744
$ws.Range("I20").Value = 6.363282705914827
$ws.Range("J20").Value = 0.1565008020490666
$ws.Range("B21").Value = 3.038547814078129
$ws.Range("C21").Value = 2.399683914589218
$ws.Range("D21").Value = 0.737773313677593
$ws.Range("E21").Value = 0.3043028472376719
$ws.Range("G21").Value = 0.002610007823128924
$ws.Range("I21").Value = 6.814115733426092
$ws.Range("J21").Value = 0.1612264703140482
$ws.Range("B22").Value = 3.113854629943432
$ws.Range("C22").Value = 2.486117961305013
$ws.Range("D22").Value = 0.7475770469571899
$ws.Range("E22").Value = 0.3091164869612086
$ws.Range("G22").Value = 0.002602825708967331
$ws.Range("I22").Value = 7.109927773486902
$ws.Range("J22").Value = 0.1644147505383984
$ws.Range("B23").Value = 3.073356501745991
$ws.Range("C23").Value = 2.43983881586621
$ws.Range("D23").Value = 0.7423066697571414
$ws.Range("E23").Value = 0.3065309439754955
$ws.Range("G23").Value = 0.002606636986712327
$ws.Range("I23").Value = 6.951932463579794
$ws.Range("J23").Value = 0.1627038505489651
$ws.Range("B24").Value = 2.926617480145126
$ws.Range("C24").Value = 2.267758770101068
$ws.Range("D24").Value = 0.7231725991360634
$ws.Range("E24").Value = 0.2970958498005274
$ws.Range("G24").Value = 0.002621556095013894
$ws.Range("I24").Value = 6.355931976645024
$ws.Range("J24").Value = 0.1564252205597114
$ws.Range("B25").Value = 2.782094804177916
$ws.Range("C25").Value = 2.088663056651626
$ws.Range("D25").Value = 0.7042570888220325
$ws.Range("E25").Value = 0.2876611729665228
$ws.Range("G25").Value = 0.002638697960490145
$ws.Range("I25").Value = 5.717254831387805
$ws.Range("J25").Value = 0.1500688710861553
